# "Separate suites for qa and stging and message update"
#
# Updates the OrderId test values on the "BOL" sheet. These order numbers
# are refreshed test fixtures (old ones expired / were consumed by the
# order-management sandbox), so the literal cell text changes while the
# existing border/fill formatting of the cells is preserved. The values
# must remain *text* (not auto-converted to numbers) since downstream
# consumers treat OrderId as a string, so we force the cell format to
# Text before writing the new value.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BOL")

# Row 2 (Parcel order): 51500899 -> 51503438
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "51503438"

# Row 3 (LTL order - Confirm): 51487922 -> 51487044
$ws.Range("A3").NumberFormat = "@"
$ws.Range("A3").Value = "51487044"

# Row 4 (LTL order - Withdraw): 51487922 -> 51487044
$ws.Range("A4").NumberFormat = "@"
$ws.Range("A4").Value = "51487044"
